$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log rows appended after the existing data (rows 2-36), covering the
# 2018.08.24 session's later timestamps.
$rows = @(
    @("2018.08.24","19:57:08","RS",10,250,0.1,0.9400000000000001,2975,0.45,0,"N/A","N/A"),
    @("2018.08.24","20:04:51","RS",10,250,0.1,0.9400000000000001,2975,0.42,0,"N/A","N/A"),
    @("2018.08.24","20:13:22","RS",10,250,0.1,0.9500000000000001,2975,0.44,0,"N/A","N/A"),
    @("2018.08.24","20:14:18","RS",10,250,0.1,0.96,2975,0.44,0,"N/A","N/A"),
    @("2018.08.24","20:22:16","RS",10,250,0.1,0.9400000000000001,2975,0.43,0,"N/A","N/A"),
    @("2018.08.24","20:22:42","RS",10,250,0.1,0.97,2975,0.39,0,"N/A","N/A"),
    @("2018.08.24","20:22:59","RS",10,250,0.1,0.92,2975,0.43,0,"N/A","N/A"),
    @("2018.08.24","20:25:44","RS",10,250,0.1,0.9500000000000001,2975,0.45,0,"N/A","N/A"),
    @("2018.08.24","20:26:09","RS",10,250,0.1,0.96,2975,0.46,0,"N/A","N/A"),
    @("2018.08.24","20:26:13","RS",10,250,0.1,0.98,2975,0.4,0,"N/A","N/A"),
    @("2018.08.24","20:29:09","RS",10,250,0.1,0.98,2975,0.4,0,"N/A","N/A"),
    @("2018.08.24","20:29:39","RS",10,250,0.1,0.99,2975,0.36,4,31.51515151515151,0.01574766423669356),
    @("2018.08.24","20:33:40","RS",10,250,0.1,0.96,2975,0.42,0,"N/A","N/A"),
    @("2018.08.24","20:33:43","RS",10,250,0.1,0.92,2975,0.42,0,"N/A","N/A"),
    @("2018.08.24","20:33:51","RS",10,250,0.1,0.96,2975,0.43,0,"N/A","N/A")
)

$startRow = 37
$endRow = $startRow + $rows.Count - 1

# Column A holds dates formatted like "2018.08.24"; force it to text first so
# Excel doesn't reinterpret the strings as date serials, then restore the
# default (unstyled) cell style once the values are in place.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($idx = 0; $idx -lt $rows.Count; $idx++) {
    $r = $startRow + $idx
    $values = $rows[$idx]

    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
    $ws.Cells.Item($r, 5).Value = $values[4]
    $ws.Cells.Item($r, 6).Value = $values[5]
    $ws.Cells.Item($r, 7).Value = $values[6]
    $ws.Cells.Item($r, 8).Value = $values[7]
    $ws.Cells.Item($r, 9).Value = $values[8]
    $ws.Cells.Item($r, 10).Value = $values[9]
    $ws.Cells.Item($r, 11).Value = $values[10]
    $ws.Cells.Item($r, 12).Value = $values[11]
}

$ws.Range("A$startRow`:A$endRow").Style = "Normal"
